$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.091.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.18%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.370.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.56%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.74%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.92%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.373.49"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.43%  "

# Row 9
$ws.Range("E9").Value = "  -0.29%  "

# Row 10
$ws.Range("E10").Value = "  +1.35%  "

# Row 11
$ws.Range("E11").Value = "  +5.62%  "

# Row 12
$ws.Range("E12").Value = "  +2.78%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.949.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.53%  "

# Row 14
$ws.Range("E14").Value = "  +0.11%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.11"
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("E16").Value = "  +4.67%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.122.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.01%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.367.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.32%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.88%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.60%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.86%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "386.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.87%  "

# Row 23
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("E24").Value = "  +2.89%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.51%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.41%  "

# Row 27
$ws.Range("E27").Value = "  +6.49%  "

# Row 28
$ws.Range("E28").Value = "  +10.24%  "

# Row 29
$ws.Range("E29").Value = "  +0.14%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.35%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.90%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.52%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.17%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.91%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.33%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.04%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0329"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.04%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0743"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.77%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.855.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.87%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.25%  "

# Row 44
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.54%  "

# Row 45
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.747"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.34%  "

# Row 46
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.49%  "

# Row 47
$ws.Range("B47").Value = "RenzoRestakedETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.414.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.58%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.45%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "299.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +13.40%  "

# Row 50
$ws.Range("E50").Value = "  -1.82%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.95%  "
